# Auto-generated Excel COM-interop script to apply F-column ("想去人数") updates
# per the commit diff, across all four worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 397
$ws.Range("F5").Value = 193
$ws.Range("F6").Value = 781
$ws.Range("F7").Value = 83
$ws.Range("F8").Value = 10067
$ws.Range("F10").Value = 3463
$ws.Range("F12").Value = 2425
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 2760
$ws.Range("F16").Value = 505
$ws.Range("F17").Value = 2135
$ws.Range("F20").Value = 16
$ws.Range("F21").Value = 381
$ws.Range("F23").Value = 124
$ws.Range("F26").Value = 203
$ws.Range("F27").Value = 611
$ws.Range("F30").Value = 1246
$ws.Range("F31").Value = 100
$ws.Range("F34").Value = 2631
$ws.Range("F35").Value = 2942
$ws.Range("F36").Value = 16
$ws.Range("F38").Value = 1023
$ws.Range("F39").Value = 385
$ws.Range("F40").Value = 5
$ws.Range("F41").Value = 1287
$ws.Range("F42").Value = 83
$ws.Range("F43").Value = 101

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 173
$ws.Range("F16").Value = 171

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 975
$ws.Range("F5").Value = 1970

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 975
$ws.Range("F5").Value = 397
$ws.Range("F8").Value = 193
$ws.Range("F9").Value = 781
$ws.Range("F10").Value = 83
$ws.Range("F11").Value = 10067
$ws.Range("F12").Value = 173
$ws.Range("F15").Value = 3463
$ws.Range("F16").Value = 2425
$ws.Range("F17").Value = 25
$ws.Range("F18").Value = 2760
$ws.Range("F20").Value = 505
$ws.Range("F21").Value = 2135
$ws.Range("F24").Value = 124
$ws.Range("F27").Value = 611
$ws.Range("F29").Value = 1246
$ws.Range("F33").Value = 2631
$ws.Range("F35").Value = 2942
$ws.Range("F36").Value = 1023
$ws.Range("F39").Value = 385
$ws.Range("F41").Value = 5
$ws.Range("F44").Value = 1287
$ws.Range("F45").Value = 83
$ws.Range("F49").Value = 171
